$d = $word.ActiveDocument

# This reverts the "cambio 3 de word" commit: the document used to be a
# single paragraph ("Cambio 2 de word") and a later commit split it into
# two paragraphs ("Cambio 2 de " + "Word" / "Cambio 3 de Word "). Undo
# that by collapsing the whole story back down to one paragraph of text,
# while leaving the "_GoBack" bookmark (anchored at the very end of the
# document) exactly where it is.

# Replace all of the document's text (everything up to, but not
# including, the final paragraph mark that the bookmark sits next to)
# with the single reverted line of text.
$body = $d.Range(0, $d.Content.End - 1)
$body.Text = "Cambio 2 de word"

# That text assignment still leaves two paragraphs behind: the new text,
# followed by an empty paragraph that holds the "_GoBack" bookmark.
# Deleting the paragraph mark that now separates them merges everything
# back into the single paragraph it used to be.
$firstPara = $d.Paragraphs(1).Range
$joinMark = $d.Range($firstPara.End - 1, $firstPara.End)
$joinMark.Delete()
